# testdata.xlsx update:
#   - DashboardData sheet renamed to CategoryData and given a header row
#     (Category Name / Category Image / Test Results) with the same bold
#     header style + row height used on the LoginData sheet.
#   - Minor column width tweaks on both sheets (cosmetic re-sizing that
#     happened while the sheets were edited in the spreadsheet app).
#   - Selection/cursor state left where the user ended up (C2:C5 on
#     LoginData, A3 on CategoryData).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LoginData")
$ws2 = $wb.Worksheets.Item("DashboardData")

# --- rename DashboardData -> CategoryData ------------------------------
$ws2.Name = "CategoryData"

# --- new header row on CategoryData ------------------------------------
$ws2.Range("A1").Value = "Category Name "
$ws2.Range("B1").Value = "Category Image"
$ws2.Range("C1").Value = "Test Results"

# match the bold header look + row height already used on LoginData row 1
$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Rows.Item(1).RowHeight = 30

# --- column width tweaks -------------------------------------------------
# Range.ColumnWidth in this host is in "characters"; the engine stores the
# sheet's <col width="..."/> as (ColumnWidth + 5/7) rounded to the nearest
# 1/7th (matching Excel's pixel-based column-width model). Back the desired
# stored width off by 5/7 so it lands on the right value.
$mdwOffset = 5.0 / 7.0

$ws2.Columns.Item(1).ColumnWidth = 27.34 - $mdwOffset
$ws2.Columns.Item(2).ColumnWidth = 36.63 - $mdwOffset
$ws2.Columns.Item(3).ColumnWidth = 25.15 - $mdwOffset

$ws1.Columns.Item(1).ColumnWidth = 22.71 - $mdwOffset
$ws1.Columns.Item(2).ColumnWidth = 25.66 - $mdwOffset
$ws1.Columns.Item(3).ColumnWidth = 47.04 - $mdwOffset
$ws1.Columns.Item(4).ColumnWidth = 7.83 - $mdwOffset

# --- leave the selection where the user finished editing -----------------
$ws2.Select()
$ws2.Range("A3").Select()

$ws1.Select()
$ws1.Range("C2:C5").Select()

Write-Output "CategoryData sheet populated and workbook updated"
